# [Kadastro App] Yeni kayit eklendi: 222 - 03.08.2025 22:07:03
# Adds a new record row to the master "Kayitlar" log sheet and to the
# matching per-district "Tarsus" log sheet.

$wb = $excel.ActiveWorkbook

$recordNo = "222"
$tarih    = "2025-08-03"
$birim    = "Tarsus"
$dosya    = "3"
$parsel   = "3"
$is       = "APL."
$personel = "NEJDET TULUKÇU (K.Teknisyeni)"

function Set-TextCell {
    param($cell, $value)
    # Force the value to be stored as text (not auto-coerced into a
    # number or date), then drop the temporary "Text" number format so
    # the cell keeps the workbook's default (unstyled) appearance.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- Sheet "Kayitlar" (master log): append new row after the last one ---
$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
$rowKayitlar = 18

Set-TextCell $wsKayitlar.Cells.Item($rowKayitlar, 1) $recordNo
Set-TextCell $wsKayitlar.Cells.Item($rowKayitlar, 2) $tarih
$wsKayitlar.Cells.Item($rowKayitlar, 3).Value = $birim
Set-TextCell $wsKayitlar.Cells.Item($rowKayitlar, 4) $dosya
Set-TextCell $wsKayitlar.Cells.Item($rowKayitlar, 5) $parsel
$wsKayitlar.Cells.Item($rowKayitlar, 6).Value = $is
$wsKayitlar.Cells.Item($rowKayitlar, 7).Value = $personel

# --- Sheet "Tarsus" (per-district log): append new row after the last one ---
$wsTarsus = $wb.Worksheets.Item("Tarsus")
$rowTarsus = 3

Set-TextCell $wsTarsus.Cells.Item($rowTarsus, 1) $recordNo
Set-TextCell $wsTarsus.Cells.Item($rowTarsus, 2) $tarih
$wsTarsus.Cells.Item($rowTarsus, 3).Value = $birim
Set-TextCell $wsTarsus.Cells.Item($rowTarsus, 4) $dosya
Set-TextCell $wsTarsus.Cells.Item($rowTarsus, 5) $parsel
$wsTarsus.Cells.Item($rowTarsus, 6).Value = $is
$wsTarsus.Cells.Item($rowTarsus, 7).Value = $personel
